# Updates LogisticRegression_optimisation results (CV timing/score columns)
# for the "Biomarkers + Clinical C / top / SHAP" sheet, plus the shared
# roc_auc per-fold string, to match the re-run optimisation values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.002811622619628906
$ws.Range("C2").Value = 0.0008941970658140667
$ws.Range("D2").Value = 0.001066207885742188
$ws.Range("E2").Value = 0.0001935308680792008
$ws.Range("R2").Value = 0.5882352941176471
$ws.Range("S2").Value = 0.6666666666666666
$ws.Range("T2").Value = 0.5454545454545454
$ws.Range("U2").Value = 0.5000000000000001
$ws.Range("V2").Value = 0.8000000000000002
$ws.Range("W2").Value = 0.6200713012477719
$ws.Range("X2").Value = 0.1053980774887464
$ws.Range("Y2").Value = 10
$ws.Range("AA2").Value = 0.805
$ws.Range("AB2").Value = 0.09994442900376631
$ws.Range("Z2").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 3
$ws.Range("B3").Value = 0.002961254119873047
$ws.Range("C3").Value = 0.0006964040200510053
$ws.Range("D3").Value = 0.00106959342956543
$ws.Range("E3").Value = 0.0003110775179630426
$ws.Range("R3").Value = 0.7058823529411764
$ws.Range("S3").Value = 0.5882352941176471
$ws.Range("T3").Value = 0.6153846153846154
$ws.Range("U3").Value = 0.5000000000000001
$ws.Range("V3").Value = 0.7272727272727272
$ws.Range("W3").Value = 0.6273549979432331
$ws.Range("X3").Value = 0.08251428349620102
$ws.Range("Y3").Value = 6
$ws.Range("AA3").Value = 0.805
$ws.Range("AB3").Value = 0.09994442900376631
$ws.Range("Z3").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 4
$ws.Range("B4").Value = 0.004684925079345703
$ws.Range("C4").Value = 0.003166224026435318
$ws.Range("D4").Value = 0.00211176872253418
$ws.Range("E4").Value = 0.0007287083051314336
$ws.Range("R4").Value = 0.5882352941176471
$ws.Range("S4").Value = 0.5
$ws.Range("T4").Value = 0.5454545454545454
$ws.Range("U4").Value = 0.5000000000000001
$ws.Range("V4").Value = 0.6153846153846154
$ws.Range("W4").Value = 0.5498148909913616
$ws.Range("X4").Value = 0.04638442729097009
$ws.Range("Y4").Value = 16
$ws.Range("AA4").Value = 0.805
$ws.Range("AB4").Value = 0.09994442900376631
$ws.Range("Z4").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 5
$ws.Range("B5").Value = 0.002970027923583985
$ws.Range("C5").Value = 0.001225764849911491
$ws.Range("D5").Value = 0.00219874382019043
$ws.Range("E5").Value = 0.00157747965582195
$ws.Range("R5").Value = 0.7499999999999999
$ws.Range("S5").Value = 0.5882352941176471
$ws.Range("T5").Value = 0.4615384615384615
$ws.Range("U5").Value = 0.3333333333333333
$ws.Range("V5").Value = 0.7272727272727272
$ws.Range("W5").Value = 0.5720759632524338
$ws.Range("X5").Value = 0.158253407555595
$ws.Range("Y5").Value = 12
$ws.Range("AA5").Value = 0.805
$ws.Range("AB5").Value = 0.09994442900376631
$ws.Range("Z5").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 6
$ws.Range("B6").Value = 0.00346522331237793
$ws.Range("C6").Value = 0.002766159851852919
$ws.Range("D6").Value = 0.001352643966674805
$ws.Range("E6").Value = 0.0007227511048941182
$ws.Range("R6").Value = 0.7058823529411764
$ws.Range("S6").Value = 0.6666666666666666
$ws.Range("T6").Value = 0.6666666666666666
$ws.Range("U6").Value = 0.4705882352941176
$ws.Range("V6").Value = 0.888888888888889
$ws.Range("W6").Value = 0.6797385620915033
$ws.Range("X6").Value = 0.1330510917647004
$ws.Range("Y6").Value = 1
$ws.Range("AA6").Value = 0.805
$ws.Range("AB6").Value = 0.09994442900376631
$ws.Range("Z6").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 7
$ws.Range("B7").Value = 0.003905487060546875
$ws.Range("C7").Value = 0.001565299971844709
$ws.Range("D7").Value = 0.001105594635009766
$ws.Range("E7").Value = 0.0005565274599280203
$ws.Range("R7").Value = 0.7058823529411764
$ws.Range("S7").Value = 0.5882352941176471
$ws.Range("T7").Value = 0.6153846153846154
$ws.Range("U7").Value = 0.5000000000000001
$ws.Range("V7").Value = 0.7272727272727272
$ws.Range("W7").Value = 0.6273549979432331
$ws.Range("X7").Value = 0.08251428349620102
$ws.Range("Y7").Value = 6
$ws.Range("AA7").Value = 0.805
$ws.Range("AB7").Value = 0.09994442900376631
$ws.Range("Z7").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 8
$ws.Range("B8").Value = 0.002144050598144531
$ws.Range("C8").Value = 0.0009256619237254871
$ws.Range("D8").Value = 0.001089811325073242
$ws.Range("E8").Value = 0.0007141400218358535
$ws.Range("R8").Value = 0.6666666666666667
$ws.Range("S8").Value = 0.7142857142857143
$ws.Range("T8").Value = 0.5454545454545454
$ws.Range("U8").Value = 0.4705882352941176
$ws.Range("V8").Value = 0.8000000000000002
$ws.Range("W8").Value = 0.6393990323402089
$ws.Range("X8").Value = 0.1178701502869808
$ws.Range("Y8").Value = 4
$ws.Range("AA8").Value = 0.805
$ws.Range("AB8").Value = 0.09994442900376631
$ws.Range("Z8").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 9
$ws.Range("B9").Value = 0.002008676528930664
$ws.Range("C9").Value = 0.0006777987020905122
$ws.Range("D9").Value = 0.0006613254547119141
$ws.Range("E9").Value = 0.0000362337885725987
$ws.Range("R9").Value = 0.7499999999999999
$ws.Range("S9").Value = 0.5882352941176471
$ws.Range("T9").Value = 0.4615384615384615
$ws.Range("U9").Value = 0.3333333333333333
$ws.Range("V9").Value = 0.7272727272727272
$ws.Range("W9").Value = 0.5720759632524338
$ws.Range("X9").Value = 0.158253407555595
$ws.Range("Y9").Value = 12
$ws.Range("AA9").Value = 0.805
$ws.Range("AB9").Value = 0.09994442900376631
$ws.Range("Z9").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 10
$ws.Range("B10").Value = 0.001474666595458984
$ws.Range("C10").Value = 0.0001740659022530774
$ws.Range("D10").Value = 0.000644063949584961
$ws.Range("E10").Value = 0.0000224967101866489
$ws.Range("R10").Value = 0.7499999999999999
$ws.Range("S10").Value = 0.5882352941176471
$ws.Range("T10").Value = 0.6666666666666666
$ws.Range("U10").Value = 0.4705882352941176
$ws.Range("V10").Value = 0.7272727272727272
$ws.Range("W10").Value = 0.6405525846702317
$ws.Range("X10").Value = 0.1017900231687662
$ws.Range("Y10").Value = 3
$ws.Range("AA10").Value = 0.805
$ws.Range("AB10").Value = 0.09994442900376631
$ws.Range("Z10").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 11
$ws.Range("B11").Value = 0.001799297332763672
$ws.Range("C11").Value = 0.00006809635310845494
$ws.Range("D11").Value = 0.0006150245666503906
$ws.Range("E11").Value = 0.000019875359839964
$ws.Range("R11").Value = 0.7058823529411764
$ws.Range("S11").Value = 0.5882352941176471
$ws.Range("T11").Value = 0.6153846153846154
$ws.Range("U11").Value = 0.5000000000000001
$ws.Range("V11").Value = 0.7272727272727272
$ws.Range("W11").Value = 0.6273549979432331
$ws.Range("X11").Value = 0.08251428349620102
$ws.Range("Y11").Value = 6
$ws.Range("AA11").Value = 0.805
$ws.Range("AB11").Value = 0.09994442900376631
$ws.Range("Z11").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 12
$ws.Range("B12").Value = 0.001299810409545898
$ws.Range("C12").Value = 0.00008699428405286021
$ws.Range("D12").Value = 0.0006105899810791016
$ws.Range("E12").Value = 0.00003204154111039021
$ws.Range("R12").Value = 0.7058823529411764
$ws.Range("S12").Value = 0.5882352941176471
$ws.Range("T12").Value = 0.6153846153846154
$ws.Range("U12").Value = 0.4705882352941176
$ws.Range("V12").Value = 0.8000000000000002
$ws.Range("W12").Value = 0.6360180995475113
$ws.Range("X12").Value = 0.1112178535385477
$ws.Range("Y12").Value = 5
$ws.Range("AA12").Value = 0.805
$ws.Range("AB12").Value = 0.09994442900376631
$ws.Range("Z12").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 13
$ws.Range("B13").Value = 0.001632308959960937
$ws.Range("C13").Value = 0.0002655784609799283
$ws.Range("D13").Value = 0.0006282806396484375
$ws.Range("E13").Value = 0.00002956098111380012
$ws.Range("R13").Value = 0.7499999999999999
$ws.Range("S13").Value = 0.5882352941176471
$ws.Range("T13").Value = 0.4615384615384615
$ws.Range("U13").Value = 0.3333333333333333
$ws.Range("V13").Value = 0.7272727272727272
$ws.Range("W13").Value = 0.5720759632524338
$ws.Range("X13").Value = 0.158253407555595
$ws.Range("Y13").Value = 12
$ws.Range("AA13").Value = 0.805
$ws.Range("AB13").Value = 0.09994442900376631
$ws.Range("Z13").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 14
$ws.Range("B14").Value = 0.002530002593994141
$ws.Range("C14").Value = 0.001488350742587671
$ws.Range("D14").Value = 0.0009104728698730469
$ws.Range("E14").Value = 0.0005524923566307804
$ws.Range("R14").Value = 0.7499999999999999
$ws.Range("S14").Value = 0.5555555555555556
$ws.Range("T14").Value = 0.6666666666666666
$ws.Range("U14").Value = 0.5333333333333333
$ws.Range("V14").Value = 0.7272727272727272
$ws.Range("W14").Value = 0.6465656565656565
$ws.Range("X14").Value = 0.08800084406566655
$ws.Range("Y14").Value = 2
$ws.Range("AA14").Value = 0.805
$ws.Range("AB14").Value = 0.09994442900376631
$ws.Range("Z14").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 15
$ws.Range("B15").Value = 0.002559614181518555
$ws.Range("C15").Value = 0.001433686733064974
$ws.Range("D15").Value = 0.0007744312286376953
$ws.Range("E15").Value = 0.0003114814195577639
$ws.Range("R15").Value = 0.7058823529411764
$ws.Range("S15").Value = 0.5882352941176471
$ws.Range("T15").Value = 0.6153846153846154
$ws.Range("U15").Value = 0.5000000000000001
$ws.Range("V15").Value = 0.7272727272727272
$ws.Range("W15").Value = 0.6273549979432331
$ws.Range("X15").Value = 0.08251428349620102
$ws.Range("Y15").Value = 6
$ws.Range("AA15").Value = 0.805
$ws.Range("AB15").Value = 0.09994442900376631
$ws.Range("Z15").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 16
$ws.Range("B16").Value = 0.001432466506958008
$ws.Range("C16").Value = 0.0001146639597491264
$ws.Range("D16").Value = 0.0005985260009765625
$ws.Range("E16").Value = 0.00001319607529059128
$ws.Range("R16").Value = 0.7058823529411764
$ws.Range("S16").Value = 0.625
$ws.Range("T16").Value = 0.4615384615384615
$ws.Range("U16").Value = 0.3076923076923077
$ws.Range("V16").Value = 0.8000000000000002
$ws.Range("W16").Value = 0.5800226244343892
$ws.Range("X16").Value = 0.1757589142019813
$ws.Range("Y16").Value = 11
$ws.Range("AA16").Value = 0.805
$ws.Range("AB16").Value = 0.09994442900376631
$ws.Range("Z16").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"

# Row 17
$ws.Range("B17").Value = 0.00157008171081543
$ws.Range("C17").Value = 0.0001842266763172697
$ws.Range("D17").Value = 0.0006055831909179688
$ws.Range("E17").Value = 0.00002400959140325437
$ws.Range("R17").Value = 0.7499999999999999
$ws.Range("S17").Value = 0.5882352941176471
$ws.Range("T17").Value = 0.4615384615384615
$ws.Range("U17").Value = 0.3333333333333333
$ws.Range("V17").Value = 0.7272727272727272
$ws.Range("W17").Value = 0.5720759632524338
$ws.Range("X17").Value = 0.158253407555595
$ws.Range("Y17").Value = 12
$ws.Range("AA17").Value = 0.805
$ws.Range("AB17").Value = 0.09994442900376631
$ws.Range("Z17").Value = "[0.95833333 0.8        0.68333333 0.71666667 0.86666667]"
